# Minor fix in TSP.
# Update the "Fitness" column (C) values for rows 2-12 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 4044
    3  = 4460
    4  = 4903
    5  = 4903
    6  = 4903
    7  = 4903
    8  = 4903
    9  = 4903
    10 = 5105
    11 = 5105
    12 = 5105
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 3).Value = $values[$row]
}
